# Refresh the cryptos list (price + 1h volume change columns), matching the
# "Updated cryptos list ... with GitHub Actions" commit. Column D (Price)
# values are assigned with a leading apostrophe so Excel stores them as text
# (matching the original inlineStr cells) instead of re-parsing e.g.
# "1.000" or "14.80" as numbers and dropping the formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''25.560.47'
$ws.Range("E2").Value = '  +2.45%  '

$ws.Range("D3").Value = '''1.664.83'
$ws.Range("E3").Value = '  +1.38%  '

$ws.Range("D4").Value = '''0.9988'
$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").Value = '''236.85'
$ws.Range("E5").Value = '  +1.72%  '

$ws.Range("D6").Value = '''1.000'
$ws.Range("E6").Value = '  -0.01%  '

$ws.Range("D7").Value = '''0.4617'
$ws.Range("E7").Value = '  -2.52%  '

$ws.Range("D8").Value = '''0.2574'
$ws.Range("E8").Value = '  -0.24%  '

$ws.Range("D9").Value = '''0.06123'
$ws.Range("E9").Value = '  +0.46%  '

$ws.Range("D10").Value = '''1.664.10'
$ws.Range("E10").Value = '  +1.27%  '

$ws.Range("D11").Value = '''0.06925'
$ws.Range("E11").Value = '  -1.49%  '

$ws.Range("D12").Value = '''14.80'
$ws.Range("E12").Value = '  +2.10%  '

$ws.Range("D13").Value = '''4.330'
$ws.Range("E13").Value = '  -0.10%  '

$ws.Range("D14").Value = '''75.06'
$ws.Range("E14").Value = '  +1.99%  '

$ws.Range("D15").Value = '''0.5722'
$ws.Range("E15").Value = '  -2.49%  '

$ws.Range("D16").Value = '''1.000'
$ws.Range("E16").Value = '  -0.02%  '

$ws.Range("E17").Value = '  -0.04%  '

$ws.Range("D18").Value = '''25.565.83'
$ws.Range("E18").Value = '  +2.43%  '

$ws.Range("D19").Value = '''0.000006671'
$ws.Range("E19").Value = '  +1.38%  '

$ws.Range("D20").Value = '''11.35'
$ws.Range("E20").Value = '  +1.37%  '

$ws.Range("D21").Value = '''1.876.42'
$ws.Range("E21").Value = '  +1.00%  '

$ws.Range("E22").Value = '  +3.22%  '

$ws.Range("D23").Value = '''8.589'
$ws.Range("E23").Value = '  +0.48%  '

$ws.Range("D24").Value = '''5.212'
$ws.Range("E24").Value = '  -0.24%  '

$ws.Range("D25").Value = '''134.15'
$ws.Range("E25").Value = '  +0.15%  '

$ws.Range("D26").Value = '''14.92'
$ws.Range("E26").Value = '  +0.23%  '

$ws.Range("D27").Value = '''1.379'
$ws.Range("E27").Value = '  -0.30%  '

$ws.Range("D28").Value = '''1.715'
$ws.Range("E28").Value = '  +5.22%  '

$ws.Range("D29").Value = '''103.95'
$ws.Range("E29").Value = '  +0.72%  '

$ws.Range("D30").Value = '''3.939'
$ws.Range("E30").Value = '  +1.69%  '

$ws.Range("D31").Value = '''0.07644'
$ws.Range("E31").Value = '  +0.79%  '

$ws.Range("D32").Value = '''3.591'
$ws.Range("E32").Value = '  +0.49%  '

$ws.Range("D33").Value = '''0.04334'
$ws.Range("E33").Value = '  +1.43%  '

$ws.Range("E34").Value = '  +1.18%  '

$ws.Range("D35").Value = '''0.6056'
$ws.Range("E35").Value = '  +2.79%  '

$ws.Range("D36").Value = '''0.9382'
$ws.Range("E36").Value = '  +1.46%  '

$ws.Range("D37").Value = '''0.9296'
$ws.Range("E37").Value = '  +6.86%  '

$ws.Range("D38").Value = '''2.432'
$ws.Range("E38").Value = '  -5.81%  '

$ws.Range("D39").Value = '''106.73'
$ws.Range("E39").Value = '  +8.19%  '

$ws.Range("D40").Value = '''0.9988'
$ws.Range("E40").Value = '  -0.12%  '

$ws.Range("D41").Value = '''1.827'
$ws.Range("E41").Value = '  +4.14%  '

$ws.Range("D42").Value = '''0.01443'
$ws.Range("E42").Value = '  -3.82%  '

$ws.Range("D43").Value = '''5.057'
$ws.Range("E43").Value = '  +8.28%  '

$ws.Range("D44").Value = '''0.3700'
$ws.Range("E44").Value = '  +0.23%  '

$ws.Range("D45").Value = '''0.1109'
$ws.Range("E45").Value = '  +0.75%  '

$ws.Range("E46").Value = '  +1.13%  '

$ws.Range("B47").Value = 'Elrond'
$ws.Range("C47").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D47").Value = '''31.20'
$ws.Range("E47").Value = '  +8.70%  '

$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").Value = '''6.068'
$ws.Range("E48").Value = '  -0.35%  '

$ws.Range("D49").Value = '''7.566'
$ws.Range("E49").Value = '  +5.83%  '

$ws.Range("E50").Value = '  +0.02%  '

$ws.Range("E51").Value = '  +0.00%  '
